$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 32
$ws.Range("H32").Value = 487
$ws.Range("I32").Value = 550.5
$ws.Range("K32").Value = 550.5
$ws.Range("M32").Value = -224.5

# ALC row 86
$ws.Range("H86").Value = 50116.145
$ws.Range("I86").Value = 57652.168
$ws.Range("J86").Value = 4900
$ws.Range("K86").Value = 57652.168
$ws.Range("L86").Value = 4900
$ws.Range("M86").Value = -56529.168
$ws.Range("N86").Value = -7146

# ALC row 89
$ws.Range("H89").Value = 50116.145
$ws.Range("I89").Value = 57652.168
$ws.Range("J89").Value = 4900
$ws.Range("K89").Value = 288260.84
$ws.Range("L89").Value = 24500
$ws.Range("M89").Value = -282644.84
$ws.Range("N89").Value = -35732

# ALC row 93
$ws.Range("H93").Value = 41333.332
$ws.Range("J93").Value = 41333.332
$ws.Range("L93").Value = 41333.332
$ws.Range("N93").Value = -46325.332

# ALC row 98
$ws.Range("H98").Value = 1679.5652
$ws.Range("I98").Value = 1220.4762
$ws.Range("J98").Value = 6500
$ws.Range("K98").Value = 1220.4762
$ws.Range("L98").Value = 6500
$ws.Range("M98").Value = 277.5237999999999
$ws.Range("N98").Value = -9496

# ALC row 112
$ws.Range("H112").Value = 6120.6665
$ws.Range("I112").Value = 23566.666
$ws.Range("J112").Value = 1759.1666
$ws.Range("K112").Value = 70699.99800000001
$ws.Range("L112").Value = 5277.4998
$ws.Range("M112").Value = -69591.99800000001
$ws.Range("N112").Value = -7493.4998

# ALC row 122
$ws.Range("H122").Value = 1679.5652
$ws.Range("I122").Value = 1220.4762
$ws.Range("J122").Value = 6500
$ws.Range("K122").Value = 3661.4286
$ws.Range("L122").Value = 19500
$ws.Range("M122").Value = -1211.4286
$ws.Range("N122").Value = -24400

# ALC row 129
$ws.Range("H129").Value = 912.72546
$ws.Range("J129").Value = 949
$ws.Range("L129").Value = 2847
$ws.Range("N129").Value = -12847

# ALC row 137
$ws.Range("H137").Value = 597797
$ws.Range("I137").Value = 2914.476
$ws.Range("J137").Value = 954726.5600000001
$ws.Range("K137").Value = 8743.428
$ws.Range("L137").Value = 2864179.68
$ws.Range("M137").Value = -6193.428
$ws.Range("N137").Value = -2869279.68

# ALC row 138
$ws.Range("H138").Value = 3572.2952
$ws.Range("I138").Value = 1918.4
$ws.Range("J138").Value = 4111.609
$ws.Range("K138").Value = 5755.200000000001
$ws.Range("L138").Value = 12334.827
$ws.Range("M138").Value = -615.2000000000007
$ws.Range("N138").Value = -22614.827

# ALC row 140
$ws.Range("H140").Value = 116956.25
$ws.Range("J140").Value = 116956.25
$ws.Range("L140").Value = 116956.25
$ws.Range("N140").Value = -127316.25

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4
$ws.Range("H4").Value = 525
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 66
$ws.Range("N4").Value = -1232

# ARM row 6
$ws.Range("H6").Value = 30000
$ws.Range("I6").Value = 30000
$ws.Range("K6").Value = 30000
$ws.Range("M6").Value = -29827

# ARM row 23
$ws.Range("H23").Value = 29999.5
$ws.Range("J23").Value = 29999.5
$ws.Range("L23").Value = 29999.5
$ws.Range("N23").Value = -30517.5

# ARM row 37
$ws.Range("H37").Value = 5034
$ws.Range("I37").Value = 5034
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 5034
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -4761
$ws.Range("N37").ClearContents()

# ARM row 44
$ws.Range("H44").Value = 57400
$ws.Range("J44").Value = 69250
$ws.Range("L44").Value = 69250
$ws.Range("N44").Value = -70226

# ARM row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# ARM row 74
$ws.Range("H74").Value = 4630.7744
$ws.Range("I74").Value = 1374.6154
$ws.Range("J74").Value = 21562.8
$ws.Range("K74").Value = 1374.6154
$ws.Range("L74").Value = 21562.8
$ws.Range("M74").Value = -500.6153999999999
$ws.Range("N74").Value = -23310.8

# ARM row 77
$ws.Range("H77").Value = 4630.7744
$ws.Range("I77").Value = 1374.6154
$ws.Range("J77").Value = 21562.8
$ws.Range("K77").Value = 6873.076999999999
$ws.Range("L77").Value = 107814
$ws.Range("M77").Value = -2505.076999999999
$ws.Range("N77").Value = -116550

# ARM row 80
$ws.Range("H80").Value = 40110
$ws.Range("J80").Value = 40110
$ws.Range("L80").Value = 40110
$ws.Range("N80").Value = -42106

# ARM row 83
$ws.Range("H83").Value = 40110
$ws.Range("J83").Value = 40110
$ws.Range("L83").Value = 120330
$ws.Range("N83").Value = -130314

$ws = $wb.Worksheets.Item("BSM")
# BSM row 134
$ws.Range("H134").Value = 45591.652
$ws.Range("I134").Value = 2208.818
$ws.Range("J134").Value = 1000014
$ws.Range("K134").Value = 6626.454000000001
$ws.Range("L134").Value = 3000042
$ws.Range("M134").Value = -4091.454000000001
$ws.Range("N134").Value = -3005112

# BSM row 140
$ws.Range("H140").Value = 54897.5
$ws.Range("J140").Value = 54897.5
$ws.Range("L140").Value = 54897.5
$ws.Range("N140").Value = -65257.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 107
$ws.Range("H107").Value = 939.88464
$ws.Range("I107").Value = 1105.4166
$ws.Range("J107").Value = 798
$ws.Range("K107").Value = 1105.4166
$ws.Range("L107").Value = 798
$ws.Range("M107").Value = 814.5834
$ws.Range("N107").Value = -4638

$ws = $wb.Worksheets.Item("CUL")
# CUL row 75
$ws.Range("H75").Value = 2471.4
$ws.Range("J75").Value = 4444
$ws.Range("L75").Value = 13332
$ws.Range("N75").Value = -15328

# CUL row 78
$ws.Range("H78").Value = 2471.4
$ws.Range("J78").Value = 4444
$ws.Range("L78").Value = 39996
$ws.Range("N78").Value = -49980

# CUL row 106
$ws.Range("H106").Value = 5888.8887
$ws.Range("J106").Value = 5888.8887
$ws.Range("L106").Value = 17666.6661
$ws.Range("N106").Value = -19558.6661

$ws = $wb.Worksheets.Item("GSM")
# GSM row 99
$ws.Range("H99").Value = 5148.4
$ws.Range("I99").Value = 5148.4
$ws.Range("K99").Value = 5148.4
$ws.Range("M99").Value = -2902.4

$ws = $wb.Worksheets.Item("LTW")
# LTW row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()

# LTW row 20
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30452

# LTW row 22
$ws.Range("H22").Value = 409.18182
$ws.Range("I22").Value = 411.875
$ws.Range("K22").Value = 411.875
$ws.Range("M22").Value = -116.875

# LTW row 27
$ws.Range("H27").Value = 409.18182
$ws.Range("I27").Value = 411.875
$ws.Range("K27").Value = 411.875
$ws.Range("M27").Value = -304.875

# LTW row 46
$ws.Range("H46").Value = 1260
$ws.Range("I46").Value = 766.6667
$ws.Range("K46").Value = 766.6667
$ws.Range("M46").Value = -578.6667

$ws = $wb.Worksheets.Item("WVR")
# WVR row 75
$ws.Range("H75").Value = 38080
$ws.Range("J75").Value = 38080
$ws.Range("L75").Value = 38080
$ws.Range("N75").Value = -39952

# WVR row 78
$ws.Range("H78").Value = 38080
$ws.Range("J78").Value = 38080
$ws.Range("L78").Value = 114240
$ws.Range("N78").Value = -123600

# WVR row 107
$ws.Range("H107").Value = 3486.1428
$ws.Range("I107").Value = 1330.3334
$ws.Range("J107").Value = 4074.0908
$ws.Range("K107").Value = 3991.0002
$ws.Range("L107").Value = 12222.2724
$ws.Range("M107").Value = -2071.0002
$ws.Range("N107").Value = -16062.2724

# WVR row 132
$ws.Range("H132").Value = 2009.7778
$ws.Range("I132").Value = 1842.1562
$ws.Range("J132").Value = 3350.75
$ws.Range("K132").Value = 5526.4686
$ws.Range("L132").Value = 10052.25
$ws.Range("M132").Value = -2996.4686
$ws.Range("N132").Value = -15112.25
